{"js": "// \"Removed mention of later use\"\n// 1) Delete the whole paragraph that tells students the exercise is\n//    needed later (\"You must finish at least exercises 1-4 incl. ...\").\n// 2) Small wording fix further down: \"ready for change? Why?\" ->\n//    \"ready for a change? Why?\"\n\nconst body = context.document.body;\n\n// --- 1. Remove the \"You must finish...\" paragraph ------------------------\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst needle = \"You must finish at least exercises 1-4 incl. of this exercise as they are used in a later exercise.\";\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(needle) !== -1) {\n    paragraphs.items[i].delete();\n    break;\n  }\n}\n\n// --- 2. \"ready for change? Why?\" -> \"ready for a change? Why?\" -----------\nconst hits = body.search(\"change? Why?\", { matchCase: true, matchWholeWord: false });\nhits.load(\"items/text\");\nawait context.sync();\n\nif (hits.items.length > 0) {\n  hits.items[0].insertText(\"a \", Word.InsertLocation.before);\n}\n\nawait context.sync();\n", "ps1": "# \"Removed mention of later use\"\n# 1) Delete the whole paragraph that tells students the exercise is\n#    needed later (\"You must finish at least exercises 1-4 incl. ...\").\n# 2) Small wording fix further down: \"ready for change? Why?\" ->\n#    \"ready for a change? Why?\"\n\n$d = $word.ActiveDocument\n\n# --- 1. Remove the \"You must finish...\" paragraph ------------------------\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*You must finish at least exercises 1-4 incl. of this exercise as they are used in a later exercise.*\") {\n        $p.Range.Delete()\n        break\n    }\n}\n\n# --- 2. \"ready for change? Why?\" -> \"ready for a change? Why?\" -----------\n$rng = $d.Content\n$find = $rng.Find\n$find.ClearFormatting()\n$find.Text = \"ready for change? Why?\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\nif ($find.Execute()) {\n    $rng.Text = \"ready for a change? Why?\"\n}\n"}
